$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of master data (Portuguese language / IST zone users)
$rows = @(
    @("por","IST","globaladmin",$true,"now()"),
    @("por","IST","service-account-mosip-resident-client",$true,"now()"),
    @("por","IST","ganesh",$true,"now()"),
    @("por","IST","officer",$true,"now()")
)

$r = 4
foreach ($row in $rows) {
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
    $r++
}

# Copy the is_active boolean cell format down to the new rows
$ws.Range("D2").Copy()
$ws.Range("D4:D7").PasteSpecial(-4122)

# Match the selection left behind in the saved workbook
$ws.Range("E2:E7").Select() | Out-Null
